$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Images for Upload")
$ws.Columns.Item(7).Delete()
